# Removing less than USD 5 price from extrapolation calibration because it
# is just a noise. This recomputes the extrapolated Risk-Neutral columns
# (ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN) for the affected rows of the
# options-calibration sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D5"  = 113498.3840899606
    "E5"  = 0.007218919331115436
    "F5"  = 0.206029345518128
    "G5"  = -0.854175765390691
    "H5"  = 7.71439265675782

    "D6"  = 114100.7878876122
    "E6"  = -0.003112932297681522
    "F6"  = 0.2361768851995746
    "G6"  = -1.314695259581175
    "H6"  = 10.27327914179687

    "D7"  = 115442.6461201255
    "E7"  = -0.008375299573469348
    "F7"  = 0.2557524574072488
    "G7"  = -1.185440944045737
    "H7"  = 7.483636105461293

    "D8"  = 115839.570541657
    "E8"  = -0.02802449743433663
    "F8"  = 0.2146091413188441
    "G8"  = -1.158970142103203
    "H8"  = 8.012766176034971

    "D9"  = 117477.4342003953
    "E9"  = -0.05307724766613111
    "F9"  = 0.312523874968035
    "G9"  = -1.473725736937873
    "H9"  = 9.686048900786094

    "D10" = 118859.4549602363
    "E10" = -0.1012159493671696
    "F10" = 0.4454089434574951
    "G10" = -1.886265077413245
    "H10" = 9.435332209698778

    "D15" = 111931.971017001
    "E15" = 0.09069927527710707
    "F15" = 0.1539168340878241
    "G15" = -0.8724840120281814
    "H15" = 9.010577648849804

    "D19" = 112084.8008748879
    "E19" = 0.04092912082759503
    "F19" = 0.1758957856638458
    "G19" = -0.5795059014213836
    "H19" = 6.44729959438628

    "D20" = 113078.0937560457
    "E20" = 0.05065033281662566
    "F20" = 0.1676756273829387
    "G20" = -0.5227636346720143
    "H20" = 5.465895759083472
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
